# TC_63769 - add "DC Unit Loading Details" mini-table (E1:E3) to both sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("Add Devices Loop A", "Update Devices")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Copy existing header/data-row formatting so the new cells pick up the
    # same styles already used by the mini lookup tables on row 7/8 (header
    # style for E1, body style for E2/E3), then fill in the values.
    $ws.Range("A7").Copy() | Out-Null
    $ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Range("A8").Copy() | Out-Null
    $ws.Range("E2:E3").PasteSpecial(-4122) | Out-Null # xlPasteFormats

    $ws.Range("E1").Value = "DC Unit Loading Details Name"
    $ws.Range("E2").Value = "Current (DC Units)"
    $ws.Range("E3").Value = "Current (worst case)"

    $ws.Range("E1:E3").Select() | Out-Null
}

$excel.CutCopyMode = 0
